# Bab 1 dan Proposal
# Appends a page break followed by a multi-level "Transaksi / Master / Report"
# bullet outline (style ListParagraph, numId=1) at the end of the document body,
# right after the "Maka dari itu ..." paragraph and before the sectPr.

$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'

# --- paragraph 1: a lone page break -----------------------------------------
$frag = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>"
$frag += "<w:pPr>" + $rPr + "</w:pPr>"
$frag += "<w:r>" + $rPr + "<w:br w:type=""page""/></w:r>"
$frag += "</w:p>"

# --- list items: level, text, and whether it carries lastRenderedPageBreak --
$items = @(
  @(0, "Transaksi", $true),
  @(1, "Transaksi Jual Beli", $false),
  @(2, "Dikerjakan oleh Novanda", $false),
  @(1, "Transaksi Restock Barang", $false),
  @(2, "Dikerjakan oleh Muhammad Amin", $false),
  @(1, "Transaksi Refund Barang", $false),
  @(2, "Dikerjakan oleh Charles", $false),
  @(0, "Master", $false),
  @(1, "Master Barang", $false),
  @(2, "Dikerjakan oleh Muhammad Amin", $false),
  @(1, "Master Transaksi", $false),
  @(2, "Dikerjakan oleh Charles", $false),
  @(1, "Master User", $false),
  @(2, "Dikerjakan oleh Novanda", $false),
  @(0, "Report", $false),
  @(1, "Report Laba Rugi", $false),
  @(2, "Dikerjakan oleh Muhammad Amin", $false),
  @(1, "Report Stock", $false),
  @(2, "Dikerjakan oleh Novanda", $false),
  @(1, "Report Penjualan", $false),
  @(2, "Dikerjakan oleh Charles", $false)
)

foreach ($it in $items) {
    $ilvl = $it[0]
    $txt = $it[1]
    $hasBreakMark = $it[2]

    $frag += "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>"
    $frag += "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""$ilvl""/><w:numId w:val=""1""/></w:numPr>" + $rPr + "</w:pPr>"
    $frag += "<w:r>" + $rPr
    if ($hasBreakMark) {
        $frag += "<w:lastRenderedPageBreak/>"
    }
    $frag += "<w:t>$txt</w:t></w:r>"
    $frag += "</w:p>"
}

$range = $d.Content
$range.Collapse(0)
$range.InsertXML($frag)

Write-Output "inserted outline"
